# Auto-generated edit script applying numeric corrections to profit-calculation
# columns (H-N) across several sheets, per the scheduled-runner update.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3088318.8
$ws.Range("I17").Value = 350
$ws.Range("J17").Value = 3335356.2
$ws.Range("K17").Value = 1050
$ws.Range("L17").Value = 10006068.6
$ws.Range("M17").Value = -882
$ws.Range("N17").Value = -10006404.6
$ws.Range("H43").Value = 8561
$ws.Range("I43").Value = 2971.6
$ws.Range("J43").Value = 11355.7
$ws.Range("K43").Value = 2971.6
$ws.Range("L43").Value = 11355.7
$ws.Range("M43").Value = -2902.6
$ws.Range("N43").Value = -11493.7
$ws.Range("H116").Value = 2754.6365
$ws.Range("I116").Value = 2515.8572
$ws.Range("K116").Value = 2515.8572
$ws.Range("M116").Value = 926.1428000000001
$ws.Range("H125").Value = 2672.4
$ws.Range("I125").Value = 5365
$ws.Range("J125").Value = 877.3333
$ws.Range("K125").Value = 48285
$ws.Range("L125").Value = 7895.9997
$ws.Range("M125").Value = -45825
$ws.Range("N125").Value = -12815.9997
$ws.Range("H129").Value = 1002.29266
$ws.Range("J129").Value = 1426
$ws.Range("L129").Value = 4278
$ws.Range("N129").Value = -14278
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 6445.273
$ws.Range("I122").Value = 6514.095
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 19542.285
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -17092.285
$ws.Range("N122").Value = -19900
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 9754.048000000001
$ws.Range("I80").Value = 22490
$ws.Range("J80").Value = 202.08333
$ws.Range("K80").Value = 22490
$ws.Range("L80").Value = 202.08333
$ws.Range("M80").Value = -21492
$ws.Range("N80").Value = -2198.08333
$ws.Range("H83").Value = 9754.048000000001
$ws.Range("I83").Value = 22490
$ws.Range("J83").Value = 202.08333
$ws.Range("K83").Value = 112450
$ws.Range("L83").Value = 1010.41665
$ws.Range("M83").Value = -107458
$ws.Range("N83").Value = -10994.41665
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4576.273
$ws.Range("I31").Value = 5673.6206
$ws.Range("J31").Value = 2454.7334
$ws.Range("K31").Value = 5673.6206
$ws.Range("L31").Value = 2454.7334
$ws.Range("M31").Value = -5378.6206
$ws.Range("N31").Value = -3044.7334
$ws.Range("H34").Value = 4576.273
$ws.Range("I34").Value = 5673.6206
$ws.Range("J34").Value = 2454.7334
$ws.Range("K34").Value = 5673.6206
$ws.Range("L34").Value = 2454.7334
$ws.Range("M34").Value = -5471.6206
$ws.Range("N34").Value = -2858.7334
$ws.Range("H58").Value = 2373.5454
$ws.Range("I58").Value = 983.6667
$ws.Range("J58").Value = 4041.4
$ws.Range("K58").Value = 983.6667
$ws.Range("L58").Value = 4041.4
$ws.Range("M58").Value = -780.6667
$ws.Range("N58").Value = -4447.4
$ws.Range("H132").Value = 2356.4285
$ws.Range("I132").Value = 2404.0454
$ws.Range("J132").Value = 2275.8462
$ws.Range("K132").Value = 7212.1362
$ws.Range("L132").Value = 6827.5386
$ws.Range("M132").Value = -4682.1362
$ws.Range("N132").Value = -11887.5386
$ws.Range("H134").Value = 1255206.9
$ws.Range("I134").Value = 2517.25
$ws.Range("J134").Value = 3402674.8
$ws.Range("K134").Value = 7551.75
$ws.Range("L134").Value = 10208024.4
$ws.Range("M134").Value = -5016.75
$ws.Range("N134").Value = -10213094.4
$ws.Range("H136").Value = 2373.5454
$ws.Range("I136").Value = 983.6667
$ws.Range("J136").Value = 4041.4
$ws.Range("K136").Value = 2951.0001
$ws.Range("L136").Value = 12124.2
$ws.Range("M136").Value = -401.0001000000002
$ws.Range("N136").Value = -17224.2
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 638.14105
$ws.Range("I113").Value = 475.193
$ws.Range("K113").Value = 1425.579
$ws.Range("M113").Value = 744.421
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3444.7354
$ws.Range("I102").Value = 3581.7097
$ws.Range("J102").Value = 2029.3334
$ws.Range("K102").Value = 3581.7097
$ws.Range("L102").Value = 2029.3334
$ws.Range("M102").Value = -1959.7097
$ws.Range("N102").Value = -5273.3334
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H113").Value = 1000011
$ws.Range("I113").Value = 1000011
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1000011
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -997841
$ws.Range("N113").ClearContents()
$ws.Range("H134").Value = 31656.5
$ws.Range("J134").Value = 31656.5
$ws.Range("L134").Value = 94969.5
$ws.Range("N134").Value = -100039.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 825.1111
$ws.Range("I46").Value = 543.9
$ws.Range("K46").Value = 543.9
$ws.Range("M46").Value = -355.9
$ws.Range("H55").Value = 317.66666
$ws.Range("I55").Value = 250.16667
$ws.Range("J55").Value = 385.16666
$ws.Range("K55").Value = 250.16667
$ws.Range("L55").Value = 385.16666
$ws.Range("M55").Value = -77.16667000000001
$ws.Range("N55").Value = -731.16666
$ws.Range("H82").Value = 1961.1428
$ws.Range("I82").Value = 1406.4615
$ws.Range("J82").Value = 2862.5
$ws.Range("K82").Value = 1406.4615
$ws.Range("L82").Value = 2862.5
$ws.Range("M82").Value = -1045.4615
$ws.Range("N82").Value = -3584.5
$ws.Range("H85").Value = 1961.1428
$ws.Range("I85").Value = 1406.4615
$ws.Range("J85").Value = 2862.5
$ws.Range("K85").Value = 1406.4615
$ws.Range("L85").Value = 2862.5
$ws.Range("M85").Value = -158.4614999999999
$ws.Range("N85").Value = -5358.5
$ws.Range("H132").Value = 14710306
$ws.Range("I132").Value = 4213.967
$ws.Range("J132").Value = 125006000
$ws.Range("K132").Value = 12641.901
$ws.Range("L132").Value = 375018000
$ws.Range("M132").Value = -10111.901
$ws.Range("N132").Value = -375023060
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2006.1111
$ws.Range("I96").Value = 1800.75
$ws.Range("J96").Value = 2170.4
$ws.Range("K96").Value = 1800.75
$ws.Range("L96").Value = 2170.4
$ws.Range("M96").Value = -427.75
$ws.Range("N96").Value = -4916.4
$ws.Range("H107").Value = 1200.1875
$ws.Range("I107").Value = 1213.5333
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 3640.5999
$ws.Range("L107").Value = 3000
$ws.Range("M107").Value = -6840
$ws.Range("N107").Value = -6840
$ws.Range("H132").Value = 2606.2856
$ws.Range("I132").Value = 2562.5
$ws.Range("J132").Value = 2664.6667
$ws.Range("K132").Value = 7687.5
$ws.Range("L132").Value = 7994.000100000001
$ws.Range("M132").Value = -5157.5
$ws.Range("N132").Value = -13054.0001
